$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Force text storage so numeric-looking strings (e.g. "191.07")
    # are not silently coerced into floating point numbers, while
    # leaving the cell style untouched (reset back to Normal after).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.845.03"
$ws.Range("E2").Value = "  +9.35%  "
$ws.Range("D3").Value = "3.511.40"
$ws.Range("E3").Value = "  +10.94%  "
$ws.Range("E4").Value = "  -0.23%  "
Set-TextValue $ws.Range("D5") "191.07"
$ws.Range("E5").Value = "  +12.21%  "
Set-TextValue $ws.Range("D6") "556.35"
$ws.Range("E6").Value = "  +9.40%  "
$ws.Range("D7").Value = "3.505.74"
$ws.Range("E7").Value = "  +11.01%  "
$ws.Range("E8").Value = "  +4.63%  "
$ws.Range("E9").Value = "  -0.19%  "
Set-TextValue $ws.Range("D10") "0.639"
$ws.Range("E10").Value = "  +8.87%  "
Set-TextValue $ws.Range("D11") "57.09"
$ws.Range("E11").Value = "  +6.34%  "
$ws.Range("E12").Value = "  +17.92%  "
$ws.Range("E13").Value = "  +11.46%  "
$ws.Range("E14").Value = "  +8.16%  "
$ws.Range("D15").Value = "4.078.56"
$ws.Range("E15").Value = "  +10.60%  "
$ws.Range("D16").Value = "3.515.08"
$ws.Range("E16").Value = "  +10.78%  "
$ws.Range("D17").Value = "68.211.98"
$ws.Range("E17").Value = "  +9.89%  "
$ws.Range("E18").Value = "  +7.45%  "
Set-TextValue $ws.Range("D19") "18.39"
$ws.Range("E19").Value = "  +9.32%  "
Set-TextValue $ws.Range("D20") "11.92"
$ws.Range("E20").Value = "  +12.05%  "
$ws.Range("E21").Value = "  +8.15%  "
Set-TextValue $ws.Range("D22") "410.12"
$ws.Range("E22").Value = "  +14.09%  "
Set-TextValue $ws.Range("D23") "3.98"
$ws.Range("E23").Value = "  +8.99%  "
Set-TextValue $ws.Range("D24") "84.75"
$ws.Range("E24").Value = "  +7.75%  "
Set-TextValue $ws.Range("D25") "11.70"
$ws.Range("E25").Value = "  +9.54%  "
Set-TextValue $ws.Range("D26") "4.22"
$ws.Range("E26").Value = "  +10.94%  "
$ws.Range("E27").Value = "  +12.64%  "
$ws.Range("E28").Value = "  +1.31%  "
Set-TextValue $ws.Range("D29") "11.91"
$ws.Range("E29").Value = "  +8.82%  "
Set-TextValue $ws.Range("D30") "8.63"
$ws.Range("E30").Value = "  +7.78%  "
Set-TextValue $ws.Range("D31") "30.61"
$ws.Range("E31").Value = "  +10.30%  "
Set-TextValue $ws.Range("D32") "682.54"
$ws.Range("E32").Value = "  +10.69%  "
$ws.Range("E33").Value = "  +8.35%  "
Set-TextValue $ws.Range("D34") "11.75"
$ws.Range("E34").Value = "  +7.86%  "
$ws.Range("E35").Value = "  +9.90%  "
Set-TextValue $ws.Range("D36") "60.50"
$ws.Range("E36").Value = "  +7.16%  "
$ws.Range("D37").Value = "0.0₃0838"
$ws.Range("E37").Value = "  +26.93%  "
Set-TextValue $ws.Range("D38") "39.14"
$ws.Range("E38").Value = "  +8.97%  "
Set-TextValue $ws.Range("D39") "0.405"
$ws.Range("E39").Value = "  +8.54%  "
$ws.Range("E40").Value = "  -0.10%  "
Set-TextValue $ws.Range("D41") "3.42"
$ws.Range("E41").Value = "  +26.74%  "
$ws.Range("E42").Value = "  +12.24%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D43") "3.04"
$ws.Range("E43").Value = "  +18.14%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D44") "2.75"
$ws.Range("E44").Value = "  +16.19%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "3.041.18"
$ws.Range("E46").Value = "  +9.11%  "
$ws.Range("E47").Value = "  +11.76%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.27"
$ws.Range("E48").Value = "  +12.10%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D49") "2.76"
$ws.Range("E49").Value = "  +5.75%  "
Set-TextValue $ws.Range("D50") "9.15"
$ws.Range("E50").Value = "  +23.46%  "
$ws.Range("E51").Value = "  +8.29%  "
